$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 547, shifting existing rows 547:568 down to 548:569.
$ws.Rows.Item(547).Insert()

# Populate the newly inserted row 547 with the new data record.
$ws.Cells.Item(547, 1).Value = 4
$ws.Cells.Item(547, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(547, 3).Value = "Los Lagos"
$ws.Cells.Item(547, 4).Value = 45075
$ws.Cells.Item(547, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(547, 5).Value = 10
$ws.Cells.Item(547, 6).Value = 100114013
$ws.Cells.Item(547, 7).Value = "Zanahoria"
$ws.Cells.Item(547, 8).Value = "Sin especificar"
$ws.Cells.Item(547, 9).Value = "Primera"
$ws.Cells.Item(547, 10).Value = 150
$ws.Cells.Item(547, 11).Value = 8000
$ws.Cells.Item(547, 12).Value = 8000
$ws.Cells.Item(547, 13).Value = 8000
$ws.Cells.Item(547, 14).Value = '$/saco 20 kilos'
$ws.Cells.Item(547, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(547, 16).Value = 400
$ws.Cells.Item(547, 17).Value = 20
$ws.Cells.Item(547, 18).Value = "Hortaliza"
